# Added ifoCAST full series evaluation: the error-triangle in the active
# sheet advances by one realized quarter. Existing rows gain one more
# realized horizon value and/or lose their now out-of-range trailing cell,
# so every populated B:K cell is rewritten in place and the cells that fall
# off the (shrinking) triangle are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (origin index 10)
$ws.Range("B2").Value = -0.4869194962242832
$ws.Range("C2").Value = -0.747338533223006
$ws.Range("D2").Value = 0.0881164629354852
$ws.Range("E2").Value = -0.01196289696713632
$ws.Range("F2").Value = -0.3657494030468326
$ws.Range("G2").Value = -0.2175720126143872
$ws.Range("H2").Value = -0.06012602361023223
$ws.Range("I2").Value = -0.560785480136303
$ws.Range("J2").Value = -0.2804275996008339
$ws.Range("K2").Value = -0.4524024574458455

# Row 3 (origin index 11)
$ws.Range("B3").Value = -0.164743996189365
$ws.Range("C3").Value = 0.6707109999691262
$ws.Range("D3").Value = 0.5706316400665047
$ws.Range("E3").Value = 0.2168451339868084
$ws.Range("F3").Value = 0.3650225244192538
$ws.Range("G3").Value = 0.5224685134234088
$ws.Range("H3").Value = 0.02180905689733798
$ws.Range("I3").Value = 0.3021669374328071
$ws.Range("J3").Value = 0.1301920795877955
$ws.Range("K3").Value = 0.436054619334127

# Row 4 (origin index 12)
$ws.Range("B4").Value = 0.5750359591597685
$ws.Range("C4").Value = 0.4749565992571469
$ws.Range("D4").Value = 0.1211700931774507
$ws.Range("E4").Value = 0.269347483609896
$ws.Range("F4").Value = 0.426793472614051
$ws.Range("G4").Value = -0.07386598391201982
$ws.Range("H4").Value = 0.2064918966234494
$ws.Range("I4").Value = 0.0345170387784377
$ws.Range("J4").Value = 0.3403795785247692
$ws.Range("K4").Value = -0.271950374162517

# Row 5 (origin index 13)
$ws.Range("B5").Value = 0.7353756362558697
$ws.Range("C5").Value = 0.3815891301761735
$ws.Range("D5").Value = 0.5297665206086188
$ws.Range("E5").Value = 0.6872125096127738
$ws.Range("F5").Value = 0.186553053086703
$ws.Range("G5").Value = 0.4669109336221722
$ws.Range("H5").Value = 0.2949360757771605
$ws.Range("I5").Value = 0.600798615523492
$ws.Range("J5").Value = -0.01153133716379418
$ws.Range("K5").Value = 0.6768400480353174

# Row 6 (origin index 14)
$ws.Range("B6").Value = -0.4538658659823178
$ws.Range("C6").Value = -0.3056884755498724
$ws.Range("D6").Value = -0.1482424865457174
$ws.Range("E6").Value = -0.6489019430717882
$ws.Range("F6").Value = -0.3685440625363191
$ws.Range("G6").Value = -0.5405189203813308
$ws.Range("H6").Value = -0.2346563806349992
$ws.Range("I6").Value = -0.8469863333222853
$ws.Range("J6").Value = -0.1586149481231739
$ws.Range("K6").Value = -0.4353061035472806

# Row 7 (origin index 15) - loses trailing K7
$ws.Range("B7").Value = -0.2056091156472509
$ws.Range("C7").Value = -0.04816312664309591
$ws.Range("D7").Value = -0.5488225831691667
$ws.Range("E7").Value = -0.2684647026336975
$ws.Range("F7").Value = -0.4404395604787092
$ws.Range("G7").Value = -0.1345770207323777
$ws.Range("H7").Value = -0.7469069734196638
$ws.Range("I7").Value = -0.05853558822055238
$ws.Range("J7").Value = -0.3352267436446591
$ws.Range("K7").ClearContents()

# Row 8 (origin index 16) - loses trailing J8
$ws.Range("B8").Value = 0.3056233794366003
$ws.Range("C8").Value = -0.1950360770894705
$ws.Range("D8").Value = 0.08532180344599868
$ws.Range("E8").Value = -0.08665305439901295
$ws.Range("F8").Value = 0.2192094853473185
$ws.Range("G8").Value = -0.3931204673399676
$ws.Range("H8").Value = 0.2952509178591439
$ws.Range("I8").Value = 0.01855976243503714
$ws.Range("J8").ClearContents()

# Row 9 (origin index 17) - loses trailing I9
$ws.Range("B9").Value = -0.3432134675219158
$ws.Range("C9").Value = -0.06285558698644665
$ws.Range("D9").Value = -0.2348304448314583
$ws.Range("E9").Value = 0.0710320949148732
$ws.Range("F9").Value = -0.541297857772413
$ws.Range("G9").Value = 0.1470735274266985
$ws.Range("H9").Value = -0.1296176279974082
$ws.Range("I9").ClearContents()

# Row 10 (origin index 18) - loses trailing H10
$ws.Range("B10").Value = -0.2203015759906016
$ws.Range("C10").Value = -0.3922764338356133
$ws.Range("D10").Value = -0.0864138940892818
$ws.Range("E10").Value = -0.698743846776568
$ws.Range("F10").Value = -0.01037246157745647
$ws.Range("G10").Value = -0.2870636170015632
$ws.Range("H10").ClearContents()

# Row 11 (origin index 19) - loses trailing G11
$ws.Range("B11").Value = 0.1083830226904575
$ws.Range("C11").Value = 0.414245562436789
$ws.Range("D11").Value = -0.1980843902504972
$ws.Range("E11").Value = 0.4902869949486143
$ws.Range("F11").Value = 0.2135958395245076
$ws.Range("G11").ClearContents()

# Row 12 (origin index 20) - loses trailing F12
$ws.Range("B12").Value = 0.1338876819013198
$ws.Range("C12").Value = -0.4784422707859664
$ws.Range("D12").Value = 0.2099291144131452
$ws.Range("E12").Value = -0.06676204101096155
$ws.Range("F12").ClearContents()

# Row 13 (origin index 21) - loses trailing E13
$ws.Range("B13").Value = -0.3064674129409547
$ws.Range("C13").Value = 0.3819039722581568
$ws.Range("D13").Value = 0.1052128168340501
$ws.Range("E13").ClearContents()

# Row 14 (origin index 22) - loses trailing D14
$ws.Range("B14").Value = 0.07604143251182532
$ws.Range("C14").Value = -0.2006497229122814
$ws.Range("D14").ClearContents()

# Row 15 (origin index 23) - loses trailing C15
$ws.Range("B15").Value = 0.4116802297750048
$ws.Range("C15").ClearContents()

# Row 16 (origin index 24) - loses trailing B16, row becomes label-only
$ws.Range("B16").ClearContents()
